$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Ev. Informe. En Relación", "Evacua Informe"),
    @("Estese al merito de autos", "Ica Informa"),
    @("Incompetencia", "Incompetencia"),
    @("Agréguese Extraordinariamente a la tabla del 05 de marzo", "Agréguese A Tabla"),
    @("Acogida sin costas", "Acoge"),
    @("Acepta competencia: Se declara admisible, 10 días. NHL ONI.", "Admisibles"),
    @("Dese Cuenta en Sala de Turno", "Dese Cuenta"),
    @("Agréguese Extraordinariamente", "Agréguese A Tabla"),
    @("Al escrito tengase presente: Previo a resolver acompañese el documento ofrecido dentro del plazo de 05 días de tener por no presentado.", "Ica Solicita Diligencia"),
    @("Inicio tramitación", "Inicio Tramitación"),
    @("Nota de Acuerdo: Redacción Valentina Salvo Oviedo", "Acuerdo"),
    @("Téngase Presente", "Téngase Presente"),
    @("Dese cuenta admisibilidad", "Dese Cuenta"),
    @("Agréguese Extraordinariamente a la tabla del 07 de marzo", "Agréguese A Tabla"),
    @("Al folio 18: Téngase presente y como se pide.", "Téngase Presente"),
    @("Al folio 19: NHL. Al segundo otrosí, como se pide.", "Nhl"),
    @("Téngase por reacusado", "Evacua Informe"),
    @("Agréguese Extraordinariamente a la tabla del 04 de marzo", "Agréguese A Tabla"),
    @("Téngase presente y por cumplido lo ordenado. En Relación", "Evacua Informe"),
    @("Rechazada sin costas", "Rechazada")
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}
